$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
    # row 33
    $ws.Range("H33").Value = 8700
    $ws.Range("I33").Value = 8700
    $ws.Range("J33").Value = 0
    $ws.Range("K33").Value = 8700
    $ws.Range("L33").Value = 0
    $ws.Range("M33").Value = -8471
    $ws.Range("N33").ClearContents()
    # row 40
    $ws.Range("H40").Value = 1880.4445
    $ws.Range("I40").Value = 1803.5714
    $ws.Range("K40").Value = 1803.5714
    $ws.Range("M40").Value = -1628.5714
    # row 99
    $ws.Range("H99").Value = 547.2222
    $ws.Range("I99").Value = 285.66666
    $ws.Range("J99").Value = 678
    $ws.Range("K99").Value = 856.9999799999999
    $ws.Range("L99").Value = 2034
    $ws.Range("M99").Value = 641.0000200000001
    $ws.Range("N99").Value = -5030
    # row 116
    $ws.Range("H116").Value = 0
    $ws.Range("I116").Value = 0
    $ws.Range("K116").Value = 0
    $ws.Range("M116").ClearContents()
    # row 138
    $ws.Range("H138").Value = 4902.1514
    $ws.Range("J138").Value = 5843.7827
    $ws.Range("L138").Value = 17531.3481
    $ws.Range("N138").Value = -27811.3481
    # row 141
    $ws.Range("H141").Value = 2904.3
    $ws.Range("I141").Value = 2904.3
    $ws.Range("K141").Value = 8712.900000000001
    $ws.Range("M141").Value = -3532.900000000001

$ws = $wb.Worksheets.Item("ARM")
    # row 26
    $ws.Range("H26").Value = 255.18182
    $ws.Range("I26").Value = 238.375
    $ws.Range("J26").Value = 300
    $ws.Range("K26").Value = 238.375
    $ws.Range("L26").Value = 300
    $ws.Range("M26").Value = 91.625
    $ws.Range("N26").Value = -960
    # row 55
    $ws.Range("H55").Value = 30000
    $ws.Range("I55").Value = 30000
    $ws.Range("K55").Value = 30000
    $ws.Range("M55").Value = -29685
    # row 63
    $ws.Range("H63").Value = 6780.4736
    $ws.Range("I63").Value = 6437.9287
    $ws.Range("K63").Value = 6437.9287
    $ws.Range("M63").Value = -5751.9287
    # row 66
    $ws.Range("H66").Value = 6780.4736
    $ws.Range("I66").Value = 6437.9287
    $ws.Range("K66").Value = 32189.6435
    $ws.Range("M66").Value = -28757.6435
    # row 97
    $ws.Range("H97").Value = 647.6
    $ws.Range("I97").Value = 334.875
    $ws.Range("K97").Value = 334.875
    $ws.Range("M97").Value = 161.125
    # row 132
    $ws.Range("H132").Value = 3064.2856
    $ws.Range("I132").Value = 2991.2307
    $ws.Range("K132").Value = 8973.6921
    $ws.Range("M132").Value = -6443.6921

$ws = $wb.Worksheets.Item("CRP")
    # row 86
    $ws.Range("H86").Value = 11499.333
    $ws.Range("I86").Value = 12999.667
    $ws.Range("J86").Value = 9999
    $ws.Range("K86").Value = 12999.667
    $ws.Range("L86").Value = 9999
    $ws.Range("M86").Value = -11876.667
    $ws.Range("N86").Value = -12245
    # row 89
    $ws.Range("H89").Value = 11499.333
    $ws.Range("I89").Value = 12999.667
    $ws.Range("J89").Value = 9999
    $ws.Range("K89").Value = 64998.335
    $ws.Range("L89").Value = 49995
    $ws.Range("M89").Value = -59382.335
    $ws.Range("N89").Value = -61227
    # row 105
    $ws.Range("H105").Value = 2866.5
    $ws.Range("I105").Value = 1570.5
    $ws.Range("K105").Value = 1570.5
    $ws.Range("M105").Value = 176.5
    # row 107
    $ws.Range("H107").Value = 1166.7894
    $ws.Range("I107").Value = 1061.0769
    $ws.Range("K107").Value = 1061.0769
    $ws.Range("M107").Value = 858.9231
    # row 132
    $ws.Range("H132").Value = 2413.647
    $ws.Range("I132").Value = 2510.3572
    $ws.Range("J132").Value = 1962.3334
    $ws.Range("K132").Value = 7531.071599999999
    $ws.Range("L132").Value = 5887.0002
    $ws.Range("M132").Value = -5001.071599999999
    $ws.Range("N132").Value = -10947.0002
    # row 134
    $ws.Range("H134").Value = 884.2632
    $ws.Range("I134").Value = 792.6429000000001
    $ws.Range("J134").Value = 1140.8
    $ws.Range("K134").Value = 2377.9287
    $ws.Range("L134").Value = 3422.4
    $ws.Range("M134").Value = 157.0712999999996
    $ws.Range("N134").Value = -8492.4

$ws = $wb.Worksheets.Item("CUL")
    # row 6
    $ws.Range("H6").Value = 518.6
    $ws.Range("J6").Value = 1492.5
    $ws.Range("L6").Value = 4477.5
    $ws.Range("N6").Value = -4703.5
    # row 10
    $ws.Range("H10").Value = 1258.8572
    $ws.Range("J10").Value = 1499.75
    $ws.Range("L10").Value = 4499.25
    $ws.Range("N10").Value = -4777.25
    # row 36
    $ws.Range("H36").Value = 16500
    $ws.Range("I36").Value = 3000
    $ws.Range("K36").Value = 9000
    $ws.Range("M36").Value = -8831
    # row 47
    $ws.Range("H47").Value = 933.8461
    $ws.Range("I47").Value = 933.8461
    $ws.Range("K47").Value = 2801.5383
    $ws.Range("M47").Value = -2370.5383
    # row 48
    $ws.Range("H48").Value = 1499
    $ws.Range("J48").Value = 1998
    $ws.Range("L48").Value = 5994
    $ws.Range("N48").Value = -6494
    # row 50
    $ws.Range("H50").Value = 945.6667
    $ws.Range("I50").Value = 987.4286
    $ws.Range("K50").Value = 2962.2858
    $ws.Range("M50").Value = -2481.2858
    # row 53
    $ws.Range("H53").Value = 945.6667
    $ws.Range("I53").Value = 987.4286
    $ws.Range("K53").Value = 2962.2858
    $ws.Range("M53").Value = -2481.2858
    # row 55
    $ws.Range("H55").Value = 2000
    $ws.Range("J55").Value = 2000
    $ws.Range("L55").Value = 6000
    $ws.Range("N55").Value = -6354
    # row 82
    $ws.Range("H82").Value = 5000
    $ws.Range("J82").Value = 5000
    $ws.Range("L82").Value = 15000
    $ws.Range("N82").Value = -15812
    # row 85
    $ws.Range("H85").Value = 5000
    $ws.Range("J85").Value = 5000
    $ws.Range("L85").Value = 15000
    $ws.Range("N85").Value = -17808
    # row 129
    $ws.Range("H129").Value = 558163.75
    $ws.Range("I129").Value = 1385.5
    $ws.Range("K129").Value = 4156.5
    $ws.Range("M129").Value = 843.5
    # row 138
    $ws.Range("H138").Value = 8334972
    $ws.Range("J138").Value = 1512.6666
    $ws.Range("L138").Value = 4537.9998
    $ws.Range("N138").Value = -14817.9998

$ws = $wb.Worksheets.Item("GSM")
    # row 82
    $ws.Range("H82").Value = 99999
    $ws.Range("I82").Value = 0
    $ws.Range("K82").Value = 0
    $ws.Range("M82").ClearContents()
    # row 85
    $ws.Range("H85").Value = 99999
    $ws.Range("I85").Value = 0
    $ws.Range("K85").Value = 0
    $ws.Range("M85").ClearContents()
    # row 102
    $ws.Range("H102").Value = 4181.5557
    $ws.Range("I102").Value = 3989.5715
    $ws.Range("J102").Value = 4853.5
    $ws.Range("K102").Value = 3989.5715
    $ws.Range("L102").Value = 4853.5
    $ws.Range("M102").Value = -2367.5715
    $ws.Range("N102").Value = -8097.5
    # row 132
    $ws.Range("H132").Value = 2487.2
    $ws.Range("I132").Value = 1696.7142
    $ws.Range("J132").Value = 4331.6665
    $ws.Range("K132").Value = 5090.142599999999
    $ws.Range("L132").Value = 12994.9995
    $ws.Range("M132").Value = -2560.142599999999
    $ws.Range("N132").Value = -18054.9995

$ws = $wb.Worksheets.Item("LTW")
    # row 122
    $ws.Range("H122").Value = 7824.778
    $ws.Range("I122").Value = 8274.929
    $ws.Range("K122").Value = 24824.787
    $ws.Range("M122").Value = -22374.787
    # row 132
    $ws.Range("H132").Value = 6033.9165
    $ws.Range("I132").Value = 6490.7
    $ws.Range("J132").Value = 3750
    $ws.Range("K132").Value = 19472.1
    $ws.Range("L132").Value = 11250
    $ws.Range("M132").Value = -16942.1
    $ws.Range("N132").Value = -16310
    # row 133
    $ws.Range("H133").Value = 0
    $ws.Range("J133").Value = 0
    $ws.Range("L133").Value = 0
    $ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
    # row 56
    $ws.Range("H56").Value = 43000
    $ws.Range("J56").Value = 43000
    $ws.Range("L56").Value = 43000
    $ws.Range("N56").Value = -44428
    # row 107
    $ws.Range("H107").Value = 722.4761999999999
    $ws.Range("I107").Value = 702.6875
    $ws.Range("J107").Value = 785.8
    $ws.Range("K107").Value = 2108.0625
    $ws.Range("L107").Value = 2357.4
    $ws.Range("M107").Value = -188.0625
    $ws.Range("N107").Value = -6197.4
    # row 113
    $ws.Range("H113").Value = 698.5
    $ws.Range("I113").Value = 622.36365
    $ws.Range("K113").Value = 1867.09095
    $ws.Range("M113").Value = 302.90905
    # row 122
    $ws.Range("H122").Value = 1641.3889
    $ws.Range("I122").Value = 1443.8235
    $ws.Range("K122").Value = 4331.470499999999
    $ws.Range("M122").Value = -1881.470499999999
